# additional screen deck updates / simplified guide page
#
# 1) Bump the cached "datetimeFigureOut" footer field from 5/6/2020 to
#    5/12/2020 everywhere it is cached (the slide master and all eleven
#    slide layouts each carry their own copy of the field's last-rendered
#    text).
# 2) Simplify the big instructional textbox on slide 5 ("Place hand over
#    top of robot ... timer runs out" -> "Place hand over robot when it
#    stops"), bump its font size way up, and re-flow/re-position the
#    textbox to its new (wider/shorter, slightly off-canvas) frame.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "5/6/2020") {
                $shp.TextFrame.TextRange.Text = "5/12/2020"
            }
        }
    }
}

# --- slide master ---
Update-DatePlaceholder $p.SlideMaster.Shapes

# --- every slide layout hanging off the master ---
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# --- slide 5: simplify the guide textbox ---
$s5 = $p.Slides.Item(5)
$guide = $s5.Shapes.Item("Rectangle 3")

# Update the run text and bump the font size first so the shape's
# auto-fit (spAutoFit) settles on its natural size *before* we stamp the
# final explicit position/size over it.
$guide.TextFrame.TextRange.Font.Size = 120
$guide.TextFrame.TextRange.Text = "Place hand over robot when it stops"

$guide.Left = -53.29520034790039
$guide.Top = 32.745750427246094
$guide.Width = 1066.59033203125
$guide.Height = 298.0828552246094

Write-Output "done"
